$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.168.06'
$ws.Range("E2").Value = '  -1.54%  '

# Row 3
$ws.Range("D3").Value = '3.584.63'
$ws.Range("E3").Value = '  -2.81%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '575.49'
$ws.Range("E5").Value = '  -5.20%  '

# Row 6
$ws.Range("D6").Value = '192.18'
$ws.Range("E6").Value = '  -0.27%  '

# Row 7
$ws.Range("D7").Value = '3.576.71'
$ws.Range("E7").Value = '  -2.77%  '

# Row 8
$ws.Range("E8").Value = '  -2.10%  '

# Row 9
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").Value = '0.679'
$ws.Range("E10").Value = '  -5.59%  '

# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '55.92'
$ws.Range("E11").Value = '  -5.97%  '

# Row 12
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '0.149'
$ws.Range("E12").Value = '  -5.76%  '

# Row 13
$ws.Range("E13").Value = '  -4.36%  '

# Row 14
$ws.Range("D14").Value = '9.86'
$ws.Range("E14").Value = '  -4.98%  '

# Row 15
$ws.Range("D15").Value = '4.153.97'
$ws.Range("E15").Value = '  -2.62%  '

# Row 16
$ws.Range("D16").Value = '3.576.78'
$ws.Range("E16").Value = '  -2.75%  '

# Row 17
$ws.Range("D17").Value = '0.126'
$ws.Range("E17").Value = '  -1.34%  '

# Row 18
$ws.Range("D18").Value = '18.38'
$ws.Range("E18").Value = '  -4.77%  '

# Row 19
$ws.Range("D19").Value = '67.091.18'
$ws.Range("E19").Value = '  -1.31%  '

# Row 20
$ws.Range("D20").Value = '12.19'
$ws.Range("E20").Value = '  -4.53%  '

# Row 21
$ws.Range("E21").Value = '  -6.63%  '

# Row 22
$ws.Range("D22").Value = '401.23'
$ws.Range("E22").Value = '  -1.21%  '

# Row 23
$ws.Range("D23").Value = '4.20'
$ws.Range("E23").Value = '  -7.93%  '

# Row 24
$ws.Range("D24").Value = '86.01'
$ws.Range("E24").Value = '  -4.11%  '

# Row 25
$ws.Range("E25").Value = '  -1.54%  '

# Row 26
$ws.Range("D26").Value = '2.94'
$ws.Range("E26").Value = '  -3.71%  '

# Row 27
$ws.Range("E27").Value = '  -3.69%  '

# Row 28
$ws.Range("D28").Value = '6.10'
$ws.Range("E28").Value = '  +1.27%  '

# Row 29
$ws.Range("D29").Value = '3.62'
$ws.Range("E29").Value = '  -3.35%  '

# Row 30
$ws.Range("E30").Value = '  -6.26%  '

# Row 31
$ws.Range("D31").Value = '7.64'
$ws.Range("E31").Value = '  +0.97%  '

# Row 32
$ws.Range("D32").Value = '31.24'
$ws.Range("E32").Value = '  -3.64%  '

# Row 33
$ws.Range("D33").Value = '636.98'
$ws.Range("E33").Value = '  +0.87%  '

# Row 34
$ws.Range("E34").Value = '  -3.99%  '

# Row 35
$ws.Range("E35").Value = '  -5.52%  '

# Row 36
$ws.Range("D36").Value = '64.04'
$ws.Range("E36").Value = '  -4.99%  '

# Row 37
$ws.Range("D37").Value = '42.32'
$ws.Range("E37").Value = '  -10.41%  '

# Row 38
$ws.Range("D38").Value = '0.399'
$ws.Range("E38").Value = '  -2.50%  '

# Row 39
$ws.Range("E39").Value = '  +0.33%  '

# Row 40
$ws.Range("E40").Value = '  -6.11%  '

# Row 41
$ws.Range("D41").Value = '3.192.28'
$ws.Range("E41").Value = '  +10.63%  '

# Row 42
$ws.Range("E42").Value = '  -3.16%  '

# Row 43
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.07%  '

# Row 44
$ws.Range("D44").Value = '2.69'
$ws.Range("E44").Value = '  +2.70%  '

# Row 45
$ws.Range("D45").Value = '2.97'
$ws.Range("E45").Value = '  -1.36%  '

# Row 46
$ws.Range("D46").Value = '0.0417'
$ws.Range("E46").Value = '  -5.62%  '

# Row 47
$ws.Range("D47").Value = '0.131'
$ws.Range("E47").Value = '  -6.08%  '

# Row 48
$ws.Range("D48").Value = '3.09'
$ws.Range("E48").Value = '  +1.28%  '

# Row 49
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '142.37'
$ws.Range("E49").Value = '  -2.42%  '

# Row 50
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.61'
$ws.Range("E50").Value = '  -2.15%  '

# Row 51
$ws.Range("D51").Value = '8.61'
$ws.Range("E51").Value = '  -6.02%  '
